$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# SW modified team 2 contacts: replace Fraserview's Team 2 captain contact
# (Julie Winton) with Cindi Ee's contact info.
$ws.Range("B17").Value = "Cindi Ee"
$ws.Range("C17").Value = "604-830-7760"
$ws.Range("D17").Value = "cindi0516@gmail.com"

$ws.Range("B16").Select()
